$wb = $excel.ActiveWorkbook

# --- Rename the existing sheet, add the new "5-sub classes" sheet right after it ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "first_experiments"

$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "5-sub classes"

# --- Copy the B2:F7 results table from first_experiments into 5-sub classes.
#     Copy cell-by-cell (rather than one bulk Range.Copy) since the source
#     range spans merged cells (C2:F2 and B4:B6) and a single bulk copy over
#     a merged region duplicates style records; per-cell copies avoid that
#     while still carrying over number formats / fonts / borders / alignment. ---
$cols = @("B", "C", "D", "E", "F")
for ($r = 2; $r -le 7; $r++) {
    foreach ($col in $cols) {
        $addr = "$col$r"
        $ws1.Range($addr).Copy($ws2.Range($addr))
    }
}

# --- Overwrite the "5-sub classes" sheet with its own experiment numbers ---
$ws2.Range("D4").Value = 84.11
$ws2.Range("E4").Value = 42.42
$ws2.Range("F4").Value = 20.59

$ws2.Range("D5").Value = 81.180000000000007
$ws2.Range("E5").Value = 41.18
$ws2.Range("F5").Value = 21.2

$ws2.Range("D6").Value = 66.03
$ws2.Range("E6").Value = 39.24
$ws2.Range("F6").Value = 19

$ws2.Range("D7").Value = 85.97
$ws2.Range("E7").Value = 40.520000000000003
$ws2.Range("F7").Value = 40

# --- Add the explanatory labels under each table ---
$ws1.Range("B10").Value = "全クラス使った"
$ws2.Range("B10").Value = "乱数で選んだ5クラスだけ使った"

# --- Restore per-sheet selections; the last Select() below leaves
#     "5-sub classes" as the active (visible) tab, matching the saved file ---
$ws1.Range("B11").Select()
$ws2.Range("L14").Select()
